$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (row 6) mirroring the existing table structure
$ws.Range("A6").Value = 12
$ws.Range("D6").Value = "type is int not string"
$ws.Range("C6").Value = 12
$ws.Range("B6").Value = "blabla"

# Move the active selection to the newly added cell, matching the saved view state
$ws.Range("D6").Select()
